$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'70.337.63"
$ws.Range("E2").Value = "  +0.15%  "
$ws.Range("D3").Value = "'3.605.68"
$ws.Range("E3").Value = "  -0.11%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'584.11"
$ws.Range("E5").Value = "  -1.25%  "
$ws.Range("D6").Value = "'190.37"
$ws.Range("E6").Value = "  +0.31%  "
$ws.Range("D7").Value = "'0.632"
$ws.Range("E7").Value = "  -1.96%  "
$ws.Range("D8").Value = "'3.598.39"
$ws.Range("E8").Value = "  -0.19%  "
$ws.Range("E9").Value = "  +0.03%  "
$ws.Range("D10").Value = "'0.182"
$ws.Range("E10").Value = "  +2.29%  "
$ws.Range("D11").Value = "'0.667"
$ws.Range("E11").Value = "  +0.57%  "
$ws.Range("D12").Value = "'56.23"
$ws.Range("E12").Value = "  -3.63%  "
$ws.Range("D13").Value = "'0.0000311"
$ws.Range("E13").Value = "  +7.67%  "
$ws.Range("D14").Value = "'9.74"
$ws.Range("E14").Value = "  -1.29%  "
$ws.Range("D15").Value = "'4.186.82"
$ws.Range("E15").Value = "  -0.08%  "
$ws.Range("D16").Value = "'20.01"
$ws.Range("E16").Value = "  +1.86%  "
$ws.Range("D17").Value = "'3.605.58"
$ws.Range("E17").Value = "  -0.18%  "
$ws.Range("D18").Value = "'70.310.41"
$ws.Range("E18").Value = "  +0.10%  "
$ws.Range("D19").Value = "'12.76"
$ws.Range("E19").Value = "  +1.28%  "
$ws.Range("E20").Value = "  +0.12%  "
$ws.Range("E21").Value = "  -0.15%  "
$ws.Range("D22").Value = "'488.23"
$ws.Range("E22").Value = "  +0.01%  "
$ws.Range("D23").Value = "'20.09"
$ws.Range("E23").Value = "  +8.08%  "
$ws.Range("D24").Value = "'4.95"
$ws.Range("E24").Value = "  -7.73%  "
$ws.Range("B25").Value = "Litecoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D25").Value = "'97.52"
$ws.Range("E25").Value = "  +7.39%  "
$ws.Range("B26").Value = "PancakeSwap"
$ws.Range("C26").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D26").Value = "'4.42"
$ws.Range("E26").Value = "  -1.22%  "
$ws.Range("D27").Value = "'3.00"
$ws.Range("E27").Value = "  -3.82%  "
$ws.Range("D28").Value = "'11.15"
$ws.Range("E28").Value = "  -0.09%  "
$ws.Range("D29").Value = "'9.59"
$ws.Range("E29").Value = "  +1.64%  "
$ws.Range("D30").Value = "'32.46"
$ws.Range("E30").Value = "  -1.07%  "
$ws.Range("D31").Value = "'7.63"
$ws.Range("E31").Value = "  -1.65%  "
$ws.Range("D32").Value = "'12.33"
$ws.Range("E32").Value = "  -0.01%  "
$ws.Range("D33").Value = "'0.120"
$ws.Range("E33").Value = "  +1.12%  "
$ws.Range("D34").Value = "'66.43"
$ws.Range("E34").Value = "  +1.19%  "
$ws.Range("D35").Value = "'581.46"
$ws.Range("E35").Value = "  -6.86%  "
$ws.Range("D36").Value = "'39.25"
$ws.Range("E36").Value = "  +2.04%  "
$ws.Range("D37").Value = "'0.0₃0819"
$ws.Range("E37").Value = "  -0.08%  "
$ws.Range("E38").Value = "  +0.07%  "
$ws.Range("E39").Value = "  -0.55%  "
$ws.Range("D40").Value = "'3.28"
$ws.Range("E40").Value = "  +20.40%  "
$ws.Range("D41").Value = "'2.90"
$ws.Range("E41").Value = "  +7.51%  "
$ws.Range("E42").Value = "  -2.97%  "
$ws.Range("D43").Value = "'0.137"
$ws.Range("E43").Value = "  -6.55%  "
$ws.Range("D44").Value = "'3.232.33"
$ws.Range("E44").Value = "  -2.46%  "
$ws.Range("E45").Value = "  -1.51%  "
$ws.Range("D46").Value = "'0.0450"
$ws.Range("E46").Value = "  -0.77%  "
$ws.Range("D47").Value = "'9.61"
$ws.Range("E47").Value = "  +5.91%  "
$ws.Range("E48").Value = "  +1.83%  "
$ws.Range("D49").Value = "'0.139"
$ws.Range("E49").Value = "  +0.42%  "
$ws.Range("E50").Value = "  -0.03%  "
$ws.Range("D51").Value = "'3.21"
$ws.Range("E51").Value = "  -3.01%  "
